$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'245.75"
$ws.Range("D2").Style = "Normal"

# Row 3
$ws.Range("D3").Value = "'25.27"
$ws.Range("D3").Style = "Normal"

# Row 4
$ws.Range("B4").Value = "LEO"
$ws.Range("C4").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D4").Value = "'3.496"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "3LEOLEO"

# Row 5
$ws.Range("B5").Value = "HuobiToken"
$ws.Range("C5").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D5").Value = "'5.050"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "4HuobiTokenHT"

# Row 6
$ws.Range("B6").Value = "Cronos"
$ws.Range("C6").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D6").Value = "'0.05597"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "5CronosCRO"

# Row 7
$ws.Range("B7").Value = "KuCoinToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("D7").Value = "'6.556"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "6KuCoinTokenKCS"

# Row 8
$ws.Range("B8").Value = "GateToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D8").Value = "'3.013"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "7GateTokenGT"

# Row 9
$ws.Range("B9").Value = "MXToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D9").Value = "'0.8166"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "8MXTokenMX"

# Row 10
$ws.Range("B10").Value = "FTXToken"
$ws.Range("C10").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D10").Value = "'0.8344"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "9FTXTokenFTT"

# Row 11
$ws.Range("B11").Value = "One"
$ws.Range("C11").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D11").Value = "'0.0005966"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "10OneONE"

# Row 12
$ws.Range("B12").Value = "WazirX"
$ws.Range("C12").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D12").Value = "'0.1335"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "11WazirXWRX"

# Row 13
$ws.Range("B13").Value = "MandalaExchangeToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D13").Value = "'0.06951"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "12MandalaExchangeTokenMDX"

# Row 14
$ws.Range("B14").Value = "BitrueCoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D14").Value = "'0.02827"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "13BitrueCoinBTR"

# Row 15
$ws.Range("B15").Value = "BitMartToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D15").Value = "'0.09387"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "14BitMartTokenBMX"

# Row 16
$ws.Range("B16").Value = "BitForexToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D16").Value = "'0.001520"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "15BitForexTokenBF"

# Row 17
$ws.Range("D17").Value = "'0.006169"
$ws.Range("D17").Style = "Normal"

# Row 18
$ws.Range("B18").Value = "BTSEToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D18").Value = "'2.092"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "17BTSETokenBTSE"

# Row 19
$ws.Range("B19").Value = "BitpandaEcosystemToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D19").Value = "'0.3188"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "18BitpandaEcosystemTokenBEST"

# Row 20
$ws.Range("B20").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C20").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D20").Value = "'0.03250"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "19LiechtensteinCryptoassetsExchangeLCX"

# Row 22
$ws.Range("D22").Value = "'3.741"
$ws.Range("D22").Style = "Normal"

# Row 23
$ws.Range("D23").Value = "'0.04693"
$ws.Range("D23").Style = "Normal"

# Row 24
$ws.Range("D24").Value = "'0.1340"
$ws.Range("D24").Style = "Normal"

# Row 25
$ws.Range("D25").Value = "'0.001242"
$ws.Range("D25").Style = "Normal"

# Row 26
$ws.Range("D26").Value = "'0.004295"
$ws.Range("D26").Style = "Normal"

# Row 27
$ws.Range("D27").Value = "'0.00009694"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "26NitroExNTX"

# Row 28
$ws.Range("D28").Value = "'0.0001939"
$ws.Range("D28").Style = "Normal"

# Row 40
$ws.Range("D40").Value = "'0.03665"
$ws.Range("D40").Style = "Normal"

# Row 41
$ws.Range("D41").Value = "'0.006187"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "40KickTokenKICKBestin24h"

# Row 42
$ws.Range("D42").Value = "'0.1053"
$ws.Range("D42").Style = "Normal"

# Row 43
$ws.Range("D43").Value = "'0.002412"
$ws.Range("D43").Style = "Normal"

# Row 44
$ws.Range("D44").Value = "'0.008215"
$ws.Range("D44").Style = "Normal"

# Row 45
$ws.Range("D45").Value = "'0.00005291"
$ws.Range("D45").Style = "Normal"

# Row 47
$ws.Range("D47").Value = "'0.1799"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "46CoinbaseStockTokenCOINWorstin24h"

# Row 48
$ws.Range("D48").Value = "'0.002015"
$ws.Range("D48").Style = "Normal"

# Row 49
$ws.Range("D49").Value = "'0.00002099"
$ws.Range("D49").Style = "Normal"

# Row 50
$ws.Range("D50").Value = "'0.0001999"
$ws.Range("D50").Style = "Normal"
